# Weekly update: insert a new record row above row 99 (shifts existing
# rows 99-105 down to 100-106) and populate it with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 99; everything below (99-105) shifts down
# one row (to 100-106), which matches the diff exactly.
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
$ws.Range("A99").Value = 6
$ws.Range("B99").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C99").Value = 'Metropolitana'
$ws.Range("D99").Value = 45142
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = 100112035
$ws.Range("G99").Value = 'Bruselas (repollito)'
$ws.Range("H99").Value = 'Sin especificar'
$ws.Range("I99").Value = 'Primera'
$ws.Range("J99").Value = 420
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 17000
$ws.Range("M99").Value = 15571
$ws.Range("N99").Value = '$/malla 15 kilos'
$ws.Range("O99").Value = 'Provincia de Quillota'
$ws.Range("P99").Value = 1038
$ws.Range("Q99").Value = 15
$ws.Range("R99").Value = 'Hortaliza'
